$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "11/15/2025"
$ws.Range("B75").Value = 0.201282096967901
$ws.Range("C75").Value = 0.798717903032099
